$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 5100
$ws.Range("E2").Value = 510
$ws.Range("F2").Value = 457
$ws.Range("G2").Value = 601
$ws.Range("H2").Value = 364
$ws.Range("I2").Value = 376
$ws.Range("J2").Value = -12
$ws.Range("K2").Value = 10482
$ws.Range("L2").Value = 2236
$ws.Range("M2").Value = 8247
$ws.Range("N2").Value = 7657
$ws.Range("O2").Value = 590
$ws.Range("P2").Value = 123
$ws.Range("Q2").Value = 70
$ws.Range("R2").Value = -281
$ws.Range("S2").Value = 19
$ws.Range("T2").Value = 246
$ws.Range("U2").Value = -176
$ws.Range("V2").Value = 1225
$ws.Range("W2").Value = 10
$ws.Range("X2").Value = 7.13
$ws.Range("Y2").Value = 5.01
$ws.Range("Z2").Value = 3.55
$ws.Range("AA2").Value = 27.11
$ws.Range("AB2").Value = 6228.44
$ws.Range("AC2").Value = 1527
$ws.Range("AD2").Value = 21.32
$ws.Range("AE2").Value = 34853
$ws.Range("AF2").Value = 0.93
$ws.Range("AG2").Value = 300
$ws.Range("AH2").Value = 0.92
$ws.Range("AI2").Value = 17.53
$ws.Range("AJ2").Value = 24630000

# Row 3
$ws.Range("D3").Value = 6168
$ws.Range("E3").Value = 661
$ws.Range("F3").Value = 661
$ws.Range("G3").Value = 696
$ws.Range("H3").Value = 729
$ws.Range("I3").Value = 744
$ws.Range("J3").Value = -15
$ws.Range("K3").Value = 9575
$ws.Range("L3").Value = 1261
$ws.Range("M3").Value = 8314
$ws.Range("N3").Value = 8314
$ws.Range("P3").Value = 123
$ws.Range("Q3").Value = 458
$ws.Range("R3").Value = -389
$ws.Range("S3").Value = -3
$ws.Range("T3").Value = 295
$ws.Range("U3").Value = 163
$ws.Range("V3").Value = 114
$ws.Range("W3").Value = 10.72
$ws.Range("X3").Value = 11.82
$ws.Range("Y3").Value = 9.31
$ws.Range("Z3").Value = 7.27
$ws.Range("AA3").Value = 15.17
$ws.Range("AB3").Value = 6778.92
$ws.Range("AC3").Value = 3020
$ws.Range("AD3").Value = 13.44
$ws.Range("AE3").Value = 37842
$ws.Range("AF3").Value = 1.07
$ws.Range("AG3").Value = 300
$ws.Range("AH3").Value = 0.74
$ws.Range("AI3").Value = 8.859999999999999
$ws.Range("AJ3").Value = 24630000
$ws.Range("O3").ClearContents()

# Row 4
$ws.Range("D4").Value = 7120
$ws.Range("E4").Value = 720
$ws.Range("F4").Value = 720
$ws.Range("G4").Value = 774
$ws.Range("H4").Value = 565
$ws.Range("I4").Value = 565
$ws.Range("K4").Value = 10380
$ws.Range("L4").Value = 1577
$ws.Range("M4").Value = 8803
$ws.Range("N4").Value = 8803
$ws.Range("P4").Value = 123
$ws.Range("Q4").Value = 180
$ws.Range("R4").Value = -443
$ws.Range("S4").Value = 219
$ws.Range("T4").Value = 475
$ws.Range("U4").Value = -295
$ws.Range("V4").Value = 384
$ws.Range("W4").Value = 10.12
$ws.Range("X4").Value = 7.93
$ws.Range("Y4").Value = 6.6
$ws.Range("Z4").Value = 5.66
$ws.Range("AA4").Value = 17.91
$ws.Range("AB4").Value = 7184.05
$ws.Range("AC4").Value = 2293
$ws.Range("AD4").Value = 15.2
$ws.Range("AE4").Value = 40071
$ws.Range("AF4").Value = 0.87
$ws.Range("AG4").Value = 300
$ws.Range("AH4").Value = 0.86
$ws.Range("AI4").Value = 11.67
$ws.Range("AJ4").Value = 24630000
$ws.Range("J4").ClearContents()
$ws.Range("O4").ClearContents()

# Row 5
$ws.Range("D5").Value = 12287
$ws.Range("E5").Value = 550
$ws.Range("F5").Value = 550
$ws.Range("G5").Value = 671
$ws.Range("H5").Value = 539
$ws.Range("I5").Value = 539
$ws.Range("K5").Value = 12523
$ws.Range("L5").Value = 3268
$ws.Range("M5").Value = 9255
$ws.Range("N5").Value = 9255
$ws.Range("P5").Value = 123
$ws.Range("Q5").Value = 728
$ws.Range("R5").Value = -1687
$ws.Range("S5").Value = 971
$ws.Range("T5").Value = 229
$ws.Range("U5").Value = 499
$ws.Range("V5").Value = 1500
$ws.Range("W5").Value = 4.48
$ws.Range("X5").Value = 4.38
$ws.Range("Y5").Value = 5.96
$ws.Range("Z5").Value = 4.7
$ws.Range("AA5").Value = 35.31
$ws.Range("AB5").Value = 7567.81
$ws.Range("AC5").Value = 2187
$ws.Range("AD5").Value = 14.38
$ws.Range("AE5").Value = 42125
$ws.Range("AF5").Value = 0.75
$ws.Range("AG5").Value = 350
$ws.Range("AH5").Value = 1.11
$ws.Range("AI5").Value = 14.28
$ws.Range("AJ5").Value = 24630000
$ws.Range("J5").ClearContents()
$ws.Range("O5").ClearContents()

# Row 6
$ws.Range("D6").Value = 12992
$ws.Range("E6").Value = 920
$ws.Range("F6").Value = 920
$ws.Range("G6").Value = 978
$ws.Range("H6").Value = 744
$ws.Range("I6").Value = 744
$ws.Range("K6").Value = 12198
$ws.Range("L6").Value = 2255
$ws.Range("M6").Value = 9943
$ws.Range("N6").Value = 9943
$ws.Range("P6").Value = 123
$ws.Range("Q6").Value = 2037
$ws.Range("R6").Value = -779
$ws.Range("S6").Value = -937
$ws.Range("T6").Value = 86
$ws.Range("U6").Value = 1951
$ws.Range("V6").Value = 648
$ws.Range("W6").Value = 7.08
$ws.Range("X6").Value = 5.72
$ws.Range("Y6").Value = 7.75
$ws.Range("Z6").Value = 6.02
$ws.Range("AA6").Value = 22.68
$ws.Range("AB6").Value = 8112.07
$ws.Range("AC6").Value = 3019
$ws.Range("AD6").Value = 12.01
$ws.Range("AE6").Value = 45259
$ws.Range("AF6").Value = 0.8
$ws.Range("AI6").Value = 11.82
$ws.Range("AJ6").Value = 24630000
$ws.Range("AG6").ClearContents()
$ws.Range("AH6").ClearContents()

# Row 7
$ws.Range("D7").Value = 12779
$ws.Range("E7").Value = 1086
$ws.Range("G7").Value = 1160
$ws.Range("H7").Value = 868
$ws.Range("I7").Value = 869
$ws.Range("K7").Value = 12779
$ws.Range("L7").Value = 2066
$ws.Range("M7").Value = 10713
$ws.Range("N7").Value = 10713
$ws.Range("P7").Value = 121
$ws.Range("Q7").Value = 1070
$ws.Range("R7").Value = -289
$ws.Range("S7").Value = -429
$ws.Range("T7").Value = 199
$ws.Range("U7").Value = 728
$ws.Range("W7").Value = 8.5
$ws.Range("X7").Value = 6.8
$ws.Range("Y7").Value = 8.41
$ws.Range("Z7").Value = 6.95
$ws.Range("AA7").Value = 19.28
$ws.Range("AC7").Value = 3527
$ws.Range("AD7").Value = 7.77
$ws.Range("AE7").Value = 48764
$ws.Range("AF7").Value = 0.5600000000000001
$ws.Range("AG7").Value = 413
$ws.Range("AH7").Value = 1.51
$ws.Range("AI7").Value = 11.7

# Row 8
$ws.Range("D8").Value = 13259
$ws.Range("E8").Value = 1185
$ws.Range("G8").Value = 1263
$ws.Range("H8").Value = 956
$ws.Range("I8").Value = 956
$ws.Range("K8").Value = 13644
$ws.Range("L8").Value = 2061
$ws.Range("M8").Value = 11582
$ws.Range("N8").Value = 11582
$ws.Range("P8").Value = 121
$ws.Range("Q8").Value = 1158
$ws.Range("R8").Value = -339
$ws.Range("S8").Value = -202
$ws.Range("T8").Value = 172
$ws.Range("U8").Value = 864
$ws.Range("W8").Value = 8.94
$ws.Range("X8").Value = 7.21
$ws.Range("Y8").Value = 8.58
$ws.Range("Z8").Value = 7.24
$ws.Range("AA8").Value = 17.8
$ws.Range("AC8").Value = 3883
$ws.Range("AD8").Value = 7.06
$ws.Range("AE8").Value = 52720
$ws.Range("AF8").Value = 0.52
$ws.Range("AG8").Value = 426
$ws.Range("AH8").Value = 1.55
$ws.Range("AI8").Value = 10.96

# Row 9
$ws.Range("D9").Value = 13741
$ws.Range("E9").Value = 1260
$ws.Range("G9").Value = 1334
$ws.Range("H9").Value = 1006
$ws.Range("I9").Value = 1006
$ws.Range("K9").Value = 14605
$ws.Range("L9").Value = 2110
$ws.Range("M9").Value = 12494
$ws.Range("N9").Value = 12494
$ws.Range("P9").Value = 121
$ws.Range("Q9").Value = 1244
$ws.Range("R9").Value = -315
$ws.Range("S9").Value = -196
$ws.Range("T9").Value = 175
$ws.Range("U9").Value = 950
$ws.Range("W9").Value = 9.17
$ws.Range("X9").Value = 7.32
$ws.Range("Y9").Value = 8.35
$ws.Range("Z9").Value = 7.12
$ws.Range("AA9").Value = 16.89
$ws.Range("AC9").Value = 4083
$ws.Range("AD9").Value = 6.71
$ws.Range("AE9").Value = 56871
$ws.Range("AF9").Value = 0.48
$ws.Range("AG9").Value = 444
$ws.Range("AH9").Value = 1.62
$ws.Range("AI9").Value = 10.89

